# Add a new worksheet "constant_names" after the last existing sheet
# (mol_ext_coefficients), populate it with the constant names "HL" and
# "H2L", select cell A2 on it, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "constant_names"

$newSheet.Range("A1").Value = "HL"
$newSheet.Range("B1").Value = "H2L"

$null = $newSheet.Range("A2").Select()
